$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.754.91"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "'2.441.26"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'567.44"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "'145.73"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.533"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("D13").Value = "'26.93"
$ws.Range("E13").Value = "  +5.63%  "
$ws.Range("D14").Value = "'0.0000181"
$ws.Range("E14").Value = "  +5.70%  "
$ws.Range("D15").Value = "'2.799.52"
$ws.Range("D16").Value = "'62.581.52"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "'2.439.01"
$ws.Range("E17").Value = "  +1.66%  "
$ws.Range("D18").Value = "'11.24"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").Value = "'324.09"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  +6.80%  "
$ws.Range("D24").Value = "'67.25"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").Value = "'8.62"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").Value = "'585.63"
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").Value = "'0.0000101"
$ws.Range("E27").Value = "  +9.06%  "
$ws.Range("D28").Value = "'2.560.23"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "'8.44"
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").Value = "'1.44"
$ws.Range("E31").Value = "  +4.36%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'1.89"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("D34").Value = "'1.53"
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("D35").Value = "'4.85"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("D38").Value = "'18.81"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("D39").Value = "'5.40"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'148.07"
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("D41").Value = "'1.81"
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'2.44"
$ws.Range("E43").Value = "  +9.20%  "
$ws.Range("D44").Value = "'148.78"
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").Value = "'0.0536"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").Value = "'20.55"
$ws.Range("E47").Value = "  +4.20%  "
$ws.Range("E48").Value = "  +2.91%  "
$ws.Range("D49").Value = "'0.0231"
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("D50").Value = "'0.0922"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("E51").Value = "  +4.49%  "
